# Apply the "Add files via upload" update to Fantasy.xlsx
# This is a weekly refresh of the fantasy-football data:
#  - Records sheet: updated win/loss/points-for numbers
#  - Schedule sheet: Week 12 games removed (season has moved on),
#    Week 13/14 projections replaced with actual results
#  - Playoffs sheet: projected-points column refreshed with actuals

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Records sheet
# ---------------------------------------------------------------------------
$records = $wb.Worksheets.Item("Records")

$recordsData = @(
    @(10, 2, 1562.6),
    @(7, 5, 1625.58),
    @(7, 5, 1654.62),
    @(9, 3, 1476.6),
    @(7, 5, 1599.6),
    @(8, 4, 1598.84),
    @(7, 5, 1617.4),
    @(5, 7, 1522.02),
    @(5, 7, 1628.54),
    @(3, 9, 1540.72),
    @(3, 9, 1491.7),
    @(1, 11, 1221.32)
)

for ($i = 0; $i -lt $recordsData.Count; $i++) {
    $r = $i + 2
    $row = $recordsData[$i]
    $records.Cells.Item($r, 2).Value = $row[0]
    $records.Cells.Item($r, 3).Value = $row[1]
    $records.Cells.Item($r, 4).Value = $row[2]
}

$records.Range("A1").Select()
$records.Range("D14").Select()

# ---------------------------------------------------------------------------
# Schedule sheet
# ---------------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

# Week 12 is no longer shown - delete its six rows (rows 2-7), shifting the
# remaining weeks up.
$schedule.Range("A2:J7").EntireRow.Delete()

# Refresh the (now-shifted) Week 13 / Week 14 rows with actual results.
$scheduleData = @(
    @(13, "Assassin's Reed", 128.1, "Death to Driscoll", 126.4),
    @(13, "One in Each Skibidi", 151.3, "The St. Brown Boy", 141.1),
    @(13, "Red Wave, Red Eyes", 145.2, "I Love Bong Pitts", 140.6),
    @(13, "Bayer Mayzyn", 137.8, "Dak White", 133),
    @(13, "ElonGPT 4o", 134.55, "London Calling", 126.2),
    @(13, "Finnegan's Fantastic Team", 125.8, "Dumpster Fire", 124),
    @(14, "Death to Driscoll", 115, "ElonGPT 4o", 133),
    @(14, "Red Wave, Red Eyes", 133.85636363636365, "One in Each Skibidi", 135),
    @(14, "The St. Brown Boy", 136, "Assassin's Reed", 124.53090909090908),
    @(14, "I Love Bong Pitts", 138, "Bayer Mayzyn", 135.21234000000001),
    @(14, "Dak White", 128, "Finnegan's Fantastic Team", 124.67843000000001),
    @(14, "London Calling", 126.2123, "Dumpster Fire", 125.78230000000001)
)

for ($i = 0; $i -lt $scheduleData.Count; $i++) {
    $r = $i + 2
    $row = $scheduleData[$i]
    $schedule.Cells.Item($r, 1).Value = $row[0]
    $schedule.Cells.Item($r, 2).Value = $row[1]
    $schedule.Cells.Item($r, 3).Value = $row[2]
    $schedule.Cells.Item($r, 4).Value = $row[3]
    $schedule.Cells.Item($r, 5).Value = $row[4]
}

$schedule.Range("E14").Select()

# ---------------------------------------------------------------------------
# Playoffs sheet
# ---------------------------------------------------------------------------
$playoffs = $wb.Worksheets.Item("Playoffs")

$playoffsData = @(141.1, 148.23419999999999, 137.80000000000001, 128.1, 134.54560000000001, 126.4, 145.19999999999999, 124, 140.6, 126.2, 133, 125.8)

for ($i = 0; $i -lt $playoffsData.Count; $i++) {
    $r = $i + 2
    $v = $playoffsData[$i]
    $playoffs.Cells.Item($r, 2).Value = $v
    $playoffs.Cells.Item($r, 3).Value = $v
    $playoffs.Cells.Item($r, 4).Value = $v
}

$playoffs.Range("F17").Select()

# ---------------------------------------------------------------------------
# Window view tweak (recorded in workbook.xml on save)
# ---------------------------------------------------------------------------
$excel.Windows.Item(1).WindowState = -4143
$schedule.Activate()
